# Edit: "Changed cleanup, updated headlines"
# Cleans up leftover <w:proofErr/> spell/grammar-check markers by merging the
# runs that Word had split around them back into single runs (text content
# is unchanged; only the run-splitting + proofErr scaffolding is removed).

$d = $word.ActiveDocument

function Set-ParagraphXmlByAnchor {
    param(
        [string]$AnchorText,
        [string]$PackageXml
    )

    $search = $d.Content
    $found = $search.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $AnchorText"
    }
    $para = $search.Paragraphs(1)
    $prange = $para.Range
    $prange.InsertXML($PackageXml)
}

# 1) Author line: "Jayanth Rao, Venkat Ramaraju" (merge off spellStart/spellEnd split)
Set-ParagraphXmlByAnchor -AnchorText "Jayanth Rao" -PackageXml '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="07B40BA1" w14:textId="77777777" w:rsidR="005941B9" w:rsidRDefault="005941B9" w:rsidP="005941B9"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Jayanth Rao, Venkat Ramaraju</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 2) "Fraudulent articles ... general public view stocks ..." paragraph
#    (merge off gramStart/gramEnd around "general public" and
#     spellStart/spellEnd around "WallStreetBets")
Set-ParagraphXmlByAnchor -AnchorText "Fraudulent articles like these" -PackageXml '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="676B3008" w14:textId="5FE1D052" w:rsidR="002E223C" w:rsidRDefault="00BE5FEB" w:rsidP="002E223C"><w:pPr><w:ind w:firstLine="270"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>With the advent of social media, the stock market has become ever more accessible to the public</w:t></w:r><w:r w:rsidR="005B286F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>, b</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">ut with this accessibility comes a risk of </w:t></w:r><w:r w:rsidR="00A92F70"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>baseless</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> opinions</w:t></w:r><w:r w:rsidR="00A92F70"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> with an amplified voice</w:t></w:r><w:r w:rsidR="005B286F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00A92F70"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="000928A6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>In 2017, the Security and Exchange Commission (SEC) announced an investigation into entities that published fraudulent articles</w:t></w:r><w:r w:rsidR="001715AA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> promoting certain stocks over others. Public companies had hired writers to publish articles without disclosing that the articles had been sponsored (SEC Press Release, 2017). </w:t></w:r><w:r w:rsidR="002E223C"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Fraudulent articles like these have the power to drastically alter how the general public view stocks and what stocks to invest in. However, there are other factors to consider as well. Online forums, as mentioned previously, have the power to create echo chambers. The same opinions can be parroted to double-down on one way of thinking. </w:t></w:r><w:r w:rsidR="00752CCF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Investment forums like Yahoo! Finance conversations and Reddit’s r/WallStreetBets can perpetuate positivity bias, where members may be </w:t></w:r><w:r w:rsidR="001F0AC9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>unnecessarily bearish or bullish</w:t></w:r><w:r w:rsidR="00752CCF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> on a given stock. Studies have shown that this </w:t></w:r><w:r w:rsidR="005956C5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">can lead to excessive trading </w:t></w:r><w:r w:rsidR="001F0AC9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>during times when that may be fiscally inadvisable</w:t></w:r><w:r w:rsidR="005956C5"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> (Tang et al., 2017). </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 3) "A large scale study to understand...:" reference paragraph
Set-ParagraphXmlByAnchor -AnchorText "A large scale study to understand" -PackageXml '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7711BADF" w14:textId="53EA3BC4" w:rsidR="005B286F" w:rsidRDefault="008C44BA" w:rsidP="005B286F"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>A large scale study to understand…:</w:t></w:r><w:r w:rsidRPr="008C44BA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>https://www.computer.org/csdl/pds/api/csdl/proceedings/download-article/12OmNzzP5HP/pdf</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# 4) "SEC.gov" reference paragraph
Set-ParagraphXmlByAnchor -AnchorText "“SEC.gov” SEC" -PackageXml '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0E5D655C" w14:textId="7618F5A2" w:rsidR="001715AA" w:rsidRDefault="001715AA" w:rsidP="005B286F"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="001715AA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>“SEC.gov” SEC: Payments for Bullish Articles on Stocks Must , 10 Apr. 2017, www.sec.gov/news/press-release/2017-79. Accessed 29 Jan. 2021.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
